$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the two new data rows for the "hidden areas" management feature
$ws.Range("A21").Value = "hidden area center"
$ws.Range("B21").Value = "y"

$ws.Range("B22").Value = "m"
$ws.Range("A22").Value = "hidden area trigger"

# Update the view so the new rows are visible, matching the authored state
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("B23").Select()
